$d = $word.ActiveDocument

# 1. Update delivery date: 01/07 -> 25/11
$d.Content.Find.Execute("01/07", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25/11", 2) | Out-Null

# 2. Split "OBS.: " so that "OBS.:" becomes bold and the trailing space
#    keeps the previous (non-bold) formatting.
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("OBS.:", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0) | Out-Null
$obsRange = $d.Range($r.Start, $r.End)
$obsRange.Font.Bold = 1

# 3. Replace the "em grupos de 4 pessoas. Como são 38 alunos, serão apenas
#    2 grupos de 3 pessoas. " sentence with the updated group-size text.
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("em grupos de 4 pessoas. Como são 38 alunos, serão apenas 2 grupos de 3 pessoas. ",
                  $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
$targetRange = $d.Range($r2.Start, $r2.End)
$targetRange.Text = "em 7 grupos de 4 pessoas (28 alunos no total). "
$targetRange.Font.Bold = 1

# 4. Move the "_GoBack" bookmark so that it sits right after "...em um"
#    (a side effect of where the last edit to this paragraph was made),
#    splitting the run that used to read "...em um arquivo ...".
$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Find.Execute("todo o código-fonte do trabalho em um", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPos) | Out-Null
